$wb = $excel.ActiveWorkbook

# Color used by the workbook's existing "HyperLink" cell style (RGB FF6495ED,
# expressed as an OLE BGR color value for Font.Color).
$hyperlinkColor = 15570276

function Set-HandbackRow {
    param($sheetName, $targetUrl, $handbackUrl, $handbackDisplay, $handbackDateTime)

    $ws = $wb.Worksheets.Item($sheetName)

    # Status: file is now handed back and in sync with en-US.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"

    # Latest Target File (E2) - same file reference as the Source File Name (A2).
    $ws.Hyperlinks.Add($ws.Range("E2"), $targetUrl, "", "", "5716bc72-aafa-4e2d-b414-735c7f57425c.md")
    $ws.Range("E2").Font.Underline = 2
    $ws.Range("E2").Font.Color = $hyperlinkColor

    # Latest Handback File (F2) - same file reference as the Latest Handoff File (C2).
    $ws.Hyperlinks.Add($ws.Range("F2"), $handbackUrl, "", "", $handbackDisplay)
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = $hyperlinkColor

    # Latest Handback DateTime (G2).
    $ws.Range("G2").Value = $handbackDateTime
}

Set-HandbackRow "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/9ef970877cebc03ab2dd7022ecbd9db425c7be05/e2e/5716bc72-aafa-4e2d-b414-735c7f57425c.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/faa79281a91dbf30985637f7d90f40fe95dd892d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5716bc72-aafa-4e2d-b414-735c7f57425c.47459c81bee3be8ae582f94547a5e79e5f7dbd89.zh-cn.xlf" `
    "5716bc72-aafa-4e2d-b414-735c7f57425c.47459c81bee3be8ae582f94547a5e79e5f7dbd89.zh-cn.xlf" `
    "2016-03-09 06:49:18"

Set-HandbackRow "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/9ef970877cebc03ab2dd7022ecbd9db425c7be05/e2e/5716bc72-aafa-4e2d-b414-735c7f57425c.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da844dbd51e189ef7c9b39855211b70775b40cb3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5716bc72-aafa-4e2d-b414-735c7f57425c.47459c81bee3be8ae582f94547a5e79e5f7dbd89.de-de.xlf" `
    "5716bc72-aafa-4e2d-b414-735c7f57425c.47459c81bee3be8ae582f94547a5e79e5f7dbd89.de-de.xlf" `
    "2016-03-09 06:49:25"
